# Updates the Price (D) and Volume(1h) (E) columns on the active sheet
# to reflect the latest cryptos snapshot, matching the source commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.349.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.633.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3819"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3566"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.220"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.414"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.287"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001229"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.633.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06950"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.562"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.31%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.360.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.550"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.37%  "
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.264"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.809.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.080"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +14.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.143"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.483"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02748"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2493"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08751"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07010"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.933"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.345"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6977"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6432"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.269"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.950"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07929"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.184"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.21%  "
